$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 18: update TX/DX/RX/IM task dates + new "no need" note ---
$ws.Range("E18").Value2 = 20200703
$ws.Range("F18").Value2 = 20200705
$ws.Range("G18").Value2 = "no need"

# --- row 22: add startdate for "class R1 > class country..." sub task ---
$ws.Range("E22").Value2 = 20200623

# --- row 25: add start/finish dates for "General trade figures, % change" ---
$ws.Range("E25").Value2 = 20200703
$ws.Range("F25").Value2 = 20200705

# --- row 26: add start/finish dates for "top X products, others..." ---
$ws.Range("E26").Value2 = 20200703
$ws.Range("F26").Value2 = 20200705

# --- move the "5 / GUI tkinter" row down from row 29 to row 34 to make
#     room for the new "excel formating" task block ---
$ws.Range("A34").Value2 = $ws.Range("A29").Value2
$ws.Range("D34").Value2 = $ws.Range("D29").Value2
$ws.Range("A29").Clear()
$ws.Range("D29").Clear()

# --- new task 6: "excel formating" / sub "c" -> "no need" ---
$ws.Range("B29").Value2 = "c"
$ws.Range("C29").Value2 = 1
$ws.Range("D28").Value2 = "excel formating"
